# Adds a new "2020-05-04" forecast-origin column (Y) and a new
# "2020-05-18" forecast-origin row (37) to both the "cases" and "deaths"
# sheets of the forecasts table, plus fills in the newly-known
# observation for 2020-05-04 (row 23 / column B).

$wb = $excel.ActiveWorkbook

# Column Y (25) values (forecast numbers) per sheet, keyed by row number.
# Rows 2-23 get an empty placeholder cell (no value yet), rows 24-36 get
# a number, and row 37 (the brand new forecast-origin row) gets its own
# single value in column Y.
$caseValuesY = @{
    24 = 34788
    25 = 36671
    26 = 38246
    27 = 39690
    28 = 40862
    29 = 41639
    30 = 42754
    31 = 43693
    32 = 44332
    33 = 45143
    34 = 45707
    35 = 46659
    36 = 47114
}
$caseY37 = 47742
$caseB23 = 32187

$deathValuesY = @{
    24 = 2856
    25 = 2990
    26 = 3097
    27 = 3194
    28 = 3270
    29 = 3314
    30 = 3390
    31 = 3452
    32 = 3489
    33 = 3543
    34 = 3577
    35 = 3644
    36 = 3671
}
$deathY37 = 3713
$deathB23 = 2654

$sheets = @("cases", "deaths")

foreach ($sheetName in $sheets) {
    $ws = $wb.Worksheets.Item($sheetName)

    if ($sheetName -eq "cases") {
        $valuesY = $caseValuesY
        $y37 = $caseY37
        $b23 = $caseB23
    } else {
        $valuesY = $deathValuesY
        $y37 = $deathY37
        $b23 = $deathB23
    }

    # --- Row 1: new header cell Y1 holding the next observation date.
    # Force text so Excel doesn't auto-convert the ISO date-looking
    # string into a date serial number (the rest of the header row is
    # plain text), then drop the temporary number format again so the
    # cell ends up with the default (unstyled) look, like its neighbours.
    $ws.Range("Y1").NumberFormat = "@"
    $ws.Range("Y1").Value = "2020-05-04"
    $ws.Range("Y1").ClearFormats()

    # --- Rows 2-23: new empty placeholder cell in column Y.
    for ($r = 2; $r -le 23; $r++) {
        $ws.Cells.Item($r, 25).ClearFormats()
    }

    # --- Rows 24-36: new forecast value in column Y.
    foreach ($r in $valuesY.Keys) {
        $ws.Cells.Item($r, 25).Value = $valuesY[$r]
    }

    # --- Row 23: the 2020-05-04 observation is now known.
    $ws.Range("B23").Value = $b23

    # --- Row 37: brand-new forecast-origin row for 2020-05-18.
    $ws.Range("A37").NumberFormat = "@"
    $ws.Range("A37").Value = "2020-05-18"
    $ws.Range("A37").ClearFormats()

    for ($c = 2; $c -le 24; $c++) {
        $ws.Cells.Item(37, $c).ClearFormats()
    }
    $ws.Cells.Item(37, 25).Value = $y37
}
